$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row
$ws.Range("D2").Value = "27.980.62"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "1.769.24"
$ws.Range("E3").Value = "  -3.31%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'321.84"
$ws.Range("E5").Value = "  -2.51%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "'0.4266"
$ws.Range("E7").Value = "  -5.48%  "
$ws.Range("D8").Value = "'0.3620"
$ws.Range("E8").Value = "  -4.81%  "
$ws.Range("D9").Value = "'43.36"
$ws.Range("E9").Value = "  -3.63%  "
$ws.Range("D10").Value = "'0.07471"
$ws.Range("E10").Value = "  -4.32%  "
$ws.Range("D11").Value = "'1.097"
$ws.Range("E11").Value = "  -3.99%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "'21.07"
$ws.Range("E13").Value = "  -5.73%  "
$ws.Range("D14").Value = "'6.109"
$ws.Range("E14").Value = "  -4.54%  "
$ws.Range("D15").Value = "'7.334"
$ws.Range("E15").Value = "  -3.10%  "
$ws.Range("D16").Value = "1.793.91"
$ws.Range("E16").Value = "  -2.42%  "
$ws.Range("D17").Value = "'92.94"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "'0.00001062"
$ws.Range("E18").Value = "  -2.56%  "
$ws.Range("D19").Value = "'0.06410"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "'17.17"
$ws.Range("E21").Value = "  -2.61%  "
$ws.Range("D22").Value = "'5.988"
$ws.Range("E22").Value = "  -6.57%  "
$ws.Range("D23").Value = "27.989.90"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").Value = "'11.34"
$ws.Range("E24").Value = "  -4.00%  "
$ws.Range("D25").Value = "'2.106"
$ws.Range("E25").Value = "  -7.90%  "
$ws.Range("D26").Value = "'158.77"
$ws.Range("E26").Value = "  +2.94%  "
$ws.Range("D27").Value = "'20.31"
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("D28").Value = "1.990.86"
$ws.Range("E28").Value = "  -2.62%  "
$ws.Range("D29").Value = "'2.165"
$ws.Range("E29").Value = "  -8.88%  "
$ws.Range("D30").Value = "'125.92"
$ws.Range("E30").Value = "  -2.92%  "
$ws.Range("D31").Value = "'1.157"
$ws.Range("E31").Value = "  -4.17%  "
$ws.Range("D32").Value = "'3.741"
$ws.Range("E32").Value = "  +2.04%  "
$ws.Range("D33").Value = "'5.642"
$ws.Range("E33").Value = "  -4.37%  "
$ws.Range("D34").Value = "'0.08914"
$ws.Range("E34").Value = "  -4.69%  "
$ws.Range("E35").Value = "  -3.57%  "
$ws.Range("D36").Value = "'0.02315"
$ws.Range("E36").Value = "  -2.17%  "
$ws.Range("D37").Value = "'0.2112"
$ws.Range("E37").Value = "  -4.42%  "
$ws.Range("D38").Value = "'5.040"
$ws.Range("E38").Value = "  -3.72%  "
$ws.Range("D41").Value = "'1.184"
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").Value = "'1.406"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "'7.838"
$ws.Range("E44").Value = "  -4.49%  "
$ws.Range("D45").Value = "'13.38"
$ws.Range("E45").Value = "  -4.49%  "
$ws.Range("D46").Value = "'0.5936"
$ws.Range("E46").Value = "  -3.89%  "
$ws.Range("D47").Value = "'3.699"
$ws.Range("E47").Value = "  -2.16%  "
$ws.Range("D48").Value = "'2.013"
$ws.Range("E48").Value = "  -2.13%  "
$ws.Range("D49").Value = "'122.80"
$ws.Range("E49").Value = "  -4.20%  "
$ws.Range("D50").Value = "'1.187"
$ws.Range("E50").Value = "  +2.61%  "
$ws.Range("D51").Value = "'0.06861"
$ws.Range("E51").Value = "  -2.26%  "

# Row 39 and 40 swap ranking positions (TheSandbox moves up, Hedera moves down) with updated values
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.6363"
$ws.Range("E39").Value = "  -5.03%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.06009"
$ws.Range("E40").Value = "  -4.95%  "

# Reset number format style on cells that were forced to text via quote-prefix,
# so no stray style/number-format is left applied to the cell.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
